$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add new row 24 data: code 99 / Altres-Diversos / Otros-Varios
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "99"
$ws.Range("B24").Value = "Altres/Diversos"
$ws.Range("C24").Value = "Otros/Varios"

# 2. Update header label "NomTaula:" -> "Nom Taula:"
$ws.Range("A1").Value = "Nom Taula:"

# 3. Expand the Excel Table (ListObject) to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:C24"))
